$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 7
$ws.Range("D11").Value = 8
$ws.Range("D16").Value = 17
$ws.Range("D21").Value = 16
$ws.Range("D26").Value = 16
$ws.Range("D36").Value = 15
$ws.Range("D41").Value = 6
$ws.Range("D51").Value = 15
$ws.Range("D56").Value = 3
$ws.Range("D80").Value = 7
$ws.Range("E80").Value = 5
$ws.Range("D81").Value = 1
$ws.Range("E81").Value = 3
$ws.Range("D86").Value = 13
$ws.Range("D91").Value = 19
$ws.Range("D96").Value = 26
$ws.Range("D101").Value = 10
$ws.Range("D116").Value = 27
$ws.Range("D121").Value = 6
$ws.Range("D126").Value = 27
$ws.Range("D131").Value = 17
$ws.Range("D141").Value = 26
$ws.Range("D146").Value = 14
$ws.Range("D150").Value = 10
$ws.Range("D151").Value = 5
$ws.Range("D171").Value = 26
$ws.Range("D176").Value = 10
$ws.Range("D200").Value = 17
$ws.Range("E200").Value = 5
$ws.Range("D201").Value = 25
$ws.Range("E201").Value = 2
$ws.Range("D206").Value = 16
$ws.Range("D209").Value = 7
$ws.Range("E209").Value = 5
$ws.Range("D210").Value = 20
$ws.Range("E210").Value = 1
$ws.Range("D211").Value = 26
$ws.Range("E211").Value = 2
$ws.Range("D230").Value = 6
$ws.Range("D231").Value = 11
$ws.Range("D241").Value = 12
$ws.Range("D256").Value = 15
$ws.Range("D265").Value = 2
$ws.Range("E265").Value = 2
$ws.Range("D266").Value = 22
$ws.Range("E266").Value = 4
$ws.Range("D271").Value = 3
$ws.Range("D285").Value = 10
$ws.Range("D291").Value = 17
$ws.Range("D296").Value = 13
$ws.Range("D306").Value = 26
$ws.Range("D331").Value = 3
$ws.Range("D356").Value = 10
$ws.Range("D366").Value = 17
$ws.Range("D375").Value = 17
$ws.Range("E375").Value = 4
$ws.Range("D376").Value = 14
$ws.Range("E376").Value = 2
$ws.Range("D391").Value = 7
$ws.Range("D396").Value = 4
$ws.Range("D401").Value = 26
$ws.Range("D406").Value = 2
$ws.Range("D421").Value = 19
$ws.Range("D426").Value = 4
$ws.Range("D431").Value = 26
$ws.Range("D436").Value = 6
$ws.Range("D441").Value = 16
$ws.Range("D451").Value = 12
$ws.Range("D466").Value = 1
$ws.Range("D471").Value = 4
$ws.Range("D476").Value = 10
$ws.Range("D486").Value = 21
$ws.Range("D496").Value = 12
$ws.Range("D505").Value = 1
$ws.Range("D506").Value = 15
$ws.Range("D510").Value = 10
$ws.Range("E510").Value = 3
$ws.Range("D511").Value = 16
$ws.Range("E511").Value = 4
$ws.Range("D516").Value = 12
$ws.Range("D520").Value = 1
$ws.Range("E520").Value = 3
$ws.Range("D521").Value = 3
$ws.Range("E521").Value = 4
$ws.Range("D541").Value = 7
$ws.Range("D556").Value = 6
$ws.Range("D561").Value = 12
$ws.Range("D575").Value = 27
$ws.Range("D576").Value = 11
$ws.Range("D581").Value = 11
$ws.Range("D586").Value = 17
$ws.Range("D591").Value = 1
$ws.Range("D616").Value = 12
$ws.Range("D621").Value = 7
$ws.Range("D626").Value = 10
$ws.Range("D636").Value = 6
$ws.Range("D646").Value = 17
$ws.Range("D651").Value = 1
$ws.Range("D656").Value = 17
$ws.Range("D661").Value = 9
$ws.Range("D666").Value = 14
$ws.Range("D667").Value = 1
$ws.Range("E667").Value = 3
$ws.Range("D668").Value = 26
$ws.Range("D669").Value = 17
$ws.Range("E669").Value = 5
$ws.Range("D670").Value = 2
$ws.Range("E670").Value = 1
$ws.Range("D671").Value = 22
$ws.Range("E671").Value = 4
$ws.Range("D672").Value = 9
$ws.Range("E672").Value = 3
$ws.Range("D673").Value = 16
$ws.Range("E673").Value = 5
$ws.Range("D674").Value = 21
$ws.Range("E674").Value = 2
$ws.Range("D675").Value = 3
$ws.Range("D676").Value = 14
